# Applies the Tab17 (DDAf_2023_Tableau_annexe_Tab17) content update:
#  - resource-rich ("*") markers move from Nigeria to Soudan du Sud / Cabo Verde
#  - the matching row shading follows the resource-rich flag
#  - the BACI/CEPII source note's "last updated" date is refreshed
#  - downstream regional aggregate figures are refreshed to match

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab17")

# --- Country-name text edits (resource-rich "*" marker moves) ------------------
$ws.Range("B34").Value = "Soudan du Sud*"   # was "Soudan du Sud"
$ws.Range("B48").Value = "Cabo Verde*"      # was "Cabo Verde"
$ws.Range("B57").Value = "Nigeria"          # was "Nigeria*"

# --- Source note: BACI/CEPII update date refreshed -----------------------------
$ws.Range("A105").Value = "Source : Calculs de l'auteur basés sur la BACI (Base de données sur le commerce international) rapportée au niveau des produits par CEPII (mise à jour le 01/02/2023)."

# --- Row shading follows the resource-rich flag ---------------------------------
# Soudan du Sud (row 34) is now resource-rich -> shade the row like the other
# shaded "*" country rows (e.g. row 17, Tchad*).
$ws.Range("B17:J17").Copy() | Out-Null
$ws.Range("B34:J34").PasteSpecial(-4122) | Out-Null

# Nigeria (row 57) is no longer resource-rich -> unshade the row like the other
# plain country rows (e.g. row 32, Seychelles).
$ws.Range("B32:J32").Copy() | Out-Null
$ws.Range("B57:J57").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# PasteSpecial(xlPasteFormats) only touches formatting, but make sure the
# displayed text is exactly as expected regardless.
$ws.Range("B34").Value = "Soudan du Sud*"
$ws.Range("B57").Value = "Nigeria"

# --- Updated aggregate figures (regional / grouping totals) --------------------
$ws.Range("C69").Value = 9418.9951000000001
$ws.Range("D69").Value = 26069.420524000001
$ws.Range("E69").Value = 2891.4281030000002
$ws.Range("F69").Value = 38379.843726999999
$ws.Range("G69").Value = 2980.318072
$ws.Range("H69").Value = 35214.188652999997
$ws.Range("I69").Value = 20001.43765
$ws.Range("J69").Value = 58195.944374999999
$ws.Range("C77").Value = 295928.30973899999
$ws.Range("D77").Value = 270682.92193999997
$ws.Range("E77").Value = 53585.117163000003
$ws.Range("F77").Value = 620196.34884200001
$ws.Range("G77").Value = 42614.062790999997
$ws.Range("H77").Value = 316869.34749100002
$ws.Range("I77").Value = 186966.86866899999
$ws.Range("J77").Value = 546450.27895099996
$ws.Range("C80").Value = 94053.410092000006
$ws.Range("D80").Value = 19858.270345000001
$ws.Range("E80").Value = 12151.048564000001
$ws.Range("F80").Value = 126062.729001
$ws.Range("G80").Value = 8819.9675220000099
$ws.Range("H80").Value = 43564.021177000002
$ws.Range("I80").Value = 25197.345923000001
$ws.Range("J80").Value = 77581.334621999995
$ws.Range("C82").Value = 157693.96486499999
$ws.Range("D82").Value = 253840.767437
$ws.Range("E82").Value = 75322.396680999998
$ws.Range("F82").Value = 486857.128983
$ws.Range("G82").Value = 62616.651850000002
$ws.Range("H82").Value = 341148.289146
$ws.Range("I82").Value = 195920.15946900001
$ws.Range("J82").Value = 599685.10046500002
$ws.Range("C84").Value = 25824.748020999999
$ws.Range("D84").Value = 70066.253211000003
$ws.Range("E84").Value = 4474.529149
$ws.Range("F84").Value = 100365.530381
$ws.Range("G84").Value = 8770.2163909999999
$ws.Range("H84").Value = 68442.571150999996
$ws.Range("I84").Value = 46829.306573000002
$ws.Range("J84").Value = 124042.094115
$ws.Range("C86").Value = 144910.870463
$ws.Range("D86").Value = 109513.34695200001
$ws.Range("E86").Value = 55869.246239
$ws.Range("F86").Value = 310293.46365400002
$ws.Range("G86").Value = 48023.400500000003
$ws.Range("H86").Value = 240343.83176299999
$ws.Range("I86").Value = 125457.55215800001
$ws.Range("J86").Value = 413824.78442099999
$ws.Range("C87").Value = 194470.135392
$ws.Range("D87").Value = 733783.27934000001
$ws.Range("E87").Value = 496735.95795900002
$ws.Range("F87").Value = 1424989.3726910001
$ws.Range("G87").Value = 297241.19151500001
$ws.Range("H87").Value = 1025532.624169
$ws.Range("I87").Value = 406270.978947
$ws.Range("J87").Value = 1729044.7946309999
$ws.Range("C89").Value = 743843.44721000001
$ws.Range("D89").Value = 2914501.4640489998
$ws.Range("E89").Value = 2353369.6726939999
$ws.Range("F89").Value = 6011714.5839529997
$ws.Range("G89").Value = 839497.87980600004
$ws.Range("H89").Value = 2355351.1080339998
$ws.Range("I89").Value = 1116544.828945
$ws.Range("J89").Value = 4311393.8167850003
$ws.Range("C90").Value = 1515861.2527719999
$ws.Range("D90").Value = 6668396.3571849996
$ws.Range("E90").Value = 4213457.9097330002
$ws.Range("F90").Value = 12397715.51969
$ws.Range("G90").Value = 1448248.0925960001
$ws.Range("H90").Value = 6581348.5886049997
$ws.Range("I90").Value = 5207118.0139340004
$ws.Range("J90").Value = 13236714.695134999
$ws.Range("C94").Value = 18286.687395000001
$ws.Range("D94").Value = 290198.90422899998
$ws.Range("E94").Value = 105489.00932300001
$ws.Range("F94").Value = 413974.60094700003
$ws.Range("G94").Value = 40252.535067999997
$ws.Range("H94").Value = 292166.44242199999
$ws.Range("I94").Value = 119966.863333
$ws.Range("J94").Value = 452385.84082300001
$ws.Range("C97").Value = 157198.79424799999
$ws.Range("D97").Value = 101755.561858
$ws.Range("E97").Value = 24326.603030999999
$ws.Range("F97").Value = 283280.95913700003
$ws.Range("G97").Value = 21273.189844
$ws.Range("H97").Value = 176202.64144499999
$ws.Range("I97").Value = 106074.482198
$ws.Range("J97").Value = 303550.31348700001
$ws.Range("C98").Value = 116682.282012
$ws.Range("D98").Value = 68947.545180000001
$ws.Range("E98").Value = 113633.696684
$ws.Range("F98").Value = 299263.52387600002
$ws.Range("G98").Value = 44715.540176000002
$ws.Range("H98").Value = 246466.85517900001
$ws.Range("I98").Value = 114928.01353700001
$ws.Range("J98").Value = 406110.40889199998
